$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 3")

# Header for new ID column
$ws.Range("F1").Value = "ID"

# ID values: rows 2-31 -> 31..60, rows 32-61 -> 31..60 (repeats)
for ($r = 2; $r -le 61; $r++) {
    $id = 31 + (($r - 2) % 30)
    $ws.Cells.Item($r, 6).Value = $id
}

# Match the final selection left behind in the sheet view
$ws.Range("F32:F61").Select()
